$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1290.3334
$ws.Range("I28").Value = 266.85715
$ws.Range("K28").Value = 266.85715
$ws.Range("M28").Value = 218.14285
$ws.Range("H107").Value = 2265.5
$ws.Range("I107").Value = 2198.9092
$ws.Range("K107").Value = 2198.9092
$ws.Range("M107").Value = -278.9092000000001
$ws.Range("H112").Value = 5151.3887
$ws.Range("J112").Value = 5366.1763
$ws.Range("L112").Value = 16098.5289
$ws.Range("N112").Value = -18314.5289
$ws.Range("H116").Value = 5297.7
$ws.Range("I116").Value = 4196.4
$ws.Range("K116").Value = 4196.4
$ws.Range("M116").Value = -754.3999999999996
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H131").Value = 1633
$ws.Range("I131").Value = 947
$ws.Range("K131").Value = 2841
$ws.Range("M131").Value = 2199

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1193.8334
$ws.Range("I4").Value = 541
$ws.Range("K4").Value = 541
$ws.Range("M4").Value = -425
$ws.Range("H5").Value = 2603.182
$ws.Range("I5").Value = 2278.3333
$ws.Range("K5").Value = 2278.3333
$ws.Range("M5").Value = -2166.3333
$ws.Range("H45").Value = 3254.3462
$ws.Range("I45").Value = 2359.2
$ws.Range("J45").Value = 4475
$ws.Range("K45").Value = 2359.2
$ws.Range("L45").Value = 4475
$ws.Range("M45").Value = -1982.2
$ws.Range("N45").Value = -5229
$ws.Range("H74").Value = 1555.3889
$ws.Range("I74").Value = 1555.3889
$ws.Range("K74").Value = 1555.3889
$ws.Range("M74").Value = -681.3888999999999
$ws.Range("H77").Value = 1555.3889
$ws.Range("I77").Value = 1555.3889
$ws.Range("K77").Value = 7776.9445
$ws.Range("M77").Value = -3408.9445
$ws.Range("H110").Value = 970.1852
$ws.Range("I110").Value = 907.8
$ws.Range("K110").Value = 907.8
$ws.Range("M110").Value = 1137.2
$ws.Range("H122").Value = 2644.7407
$ws.Range("I122").Value = 2495.92
$ws.Range("K122").Value = 7487.76
$ws.Range("M122").Value = -5037.76

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2603.182
$ws.Range("I4").Value = 2278.3333
$ws.Range("K4").Value = 2278.3333
$ws.Range("M4").Value = -2163.3333
$ws.Range("H20").Value = 3597.5
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 705.5
$ws.Range("I22").Value = 705.5
$ws.Range("K22").Value = 705.5
$ws.Range("M22").Value = -532.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12458.143
$ws.Range("I16").Value = 958
$ws.Range("K16").Value = 958
$ws.Range("M16").Value = -671
$ws.Range("H22").Value = 513.5
$ws.Range("I22").Value = 532.2308
$ws.Range("J22").Value = 432.33334
$ws.Range("K22").Value = 532.2308
$ws.Range("L22").Value = 432.33334
$ws.Range("M22").Value = -182.2308
$ws.Range("N22").Value = -1132.33334
$ws.Range("H31").Value = 3737.5625
$ws.Range("I31").Value = 4016.1738
$ws.Range("J31").Value = 3025.5557
$ws.Range("K31").Value = 4016.1738
$ws.Range("L31").Value = 3025.5557
$ws.Range("M31").Value = -3721.1738
$ws.Range("N31").Value = -3615.5557
$ws.Range("H34").Value = 3737.5625
$ws.Range("I34").Value = 4016.1738
$ws.Range("J34").Value = 3025.5557
$ws.Range("K34").Value = 4016.1738
$ws.Range("L34").Value = 3025.5557
$ws.Range("M34").Value = -3814.1738
$ws.Range("N34").Value = -3429.5557
$ws.Range("H41").Value = 2875
$ws.Range("I41").Value = 2875
$ws.Range("K41").Value = 2875
$ws.Range("M41").Value = -2447
$ws.Range("H47").Value = 29492.5
$ws.Range("I47").Value = 24950
$ws.Range("K47").Value = 24950
$ws.Range("M47").Value = -24384
$ws.Range("H113").Value = 12458.143
$ws.Range("I113").Value = 958
$ws.Range("K113").Value = 958
$ws.Range("M113").Value = 1212
$ws.Range("H132").Value = 4097.2354
$ws.Range("I132").Value = 3977
$ws.Range("J132").Value = 4658.3335
$ws.Range("K132").Value = 11931
$ws.Range("L132").Value = 13975.0005
$ws.Range("M132").Value = -9401
$ws.Range("N132").Value = -19035.0005
$ws.Range("H134").Value = 3941.2942
$ws.Range("I134").Value = 3839.0322
$ws.Range("K134").Value = 11517.0966
$ws.Range("M134").Value = -8982.096600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1253.9286
$ws.Range("I2").Value = 1449.909
$ws.Range("K2").Value = 8699.454000000002
$ws.Range("M2").Value = -8586.454000000002
$ws.Range("H115").Value = 2964
$ws.Range("I115").Value = 2964
$ws.Range("K115").Value = 8892
$ws.Range("M115").Value = -7717
$ws.Range("H131").Value = 1732.5264
$ws.Range("J131").Value = 2295.0435
$ws.Range("L131").Value = 6885.130500000001
$ws.Range("N131").Value = -16965.1305

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H102").Value = 2452.6667
$ws.Range("I102").Value = 2452.6667
$ws.Range("K102").Value = 2452.6667
$ws.Range("M102").Value = -830.6667000000002
$ws.Range("H122").Value = 3444.2693
$ws.Range("I122").Value = 3080.182
$ws.Range("K122").Value = 9240.545999999998
$ws.Range("M122").Value = -6790.545999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5332.6665
$ws.Range("I61").Value = 5499
$ws.Range("K61").Value = 5499
$ws.Range("M61").Value = -5297
$ws.Range("H100").Value = 5600
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H113").Value = 5332.6665
$ws.Range("I113").Value = 5499
$ws.Range("K113").Value = 5499
$ws.Range("M113").Value = -3329
$ws.Range("H132").Value = 2910.459
$ws.Range("I132").Value = 2933.4897
$ws.Range("J132").Value = 2816.4167
$ws.Range("K132").Value = 8800.4691
$ws.Range("L132").Value = 8449.250100000001
$ws.Range("M132").Value = -6270.4691
$ws.Range("N132").Value = -13509.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 29825.4
$ws.Range("I45").Value = 45951.5
$ws.Range("K45").Value = 45951.5
$ws.Range("M45").Value = -45460.5
$ws.Range("H113").Value = 1025.7142
$ws.Range("I113").Value = 1025.7142
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3077.1426
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -907.1425999999997
$ws.Range("N113").ClearContents()
